# The deck currently renders slides with the "Integral" theme palette
# (stored in ppt/theme/theme2.xml, the theme referenced by the one and
# only slide master) while ppt/theme/theme1.xml - reachable only from
# the notes master - holds the stock "Office Theme" palette.
#
# The authored commit swaps the contents of theme1.xml and theme2.xml
# (file names / relationships are untouched): the slide-facing theme
# becomes the stock "Office Theme" colours, and the notes-only theme
# becomes the "Integral" colours.
#
# PowerPoint's object model exposes the slide-facing theme's 12-colour
# DrawingML colour scheme via Slide.ThemeColorScheme (index order:
# 1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6 11=hlink 12=folHlink), each
# entry's .RGB being a normal COM BGR-packed RGB() value, so recolour
# every slide-facing slot to the "Office Theme" values.

function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the stock "Office Theme" colour scheme, in
# a:clrScheme slot order (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = HexToComRgb $officeTheme[$i - 1]
}
